# Refresh Volume(1h) % and Price figures for the cryptos list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.374.11"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "2.430.84"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'570.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "'143.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").Value = "2.426.11"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("E10").Value = "  -4.32%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").Value = "'5.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "'26.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("E15").Value = "  -4.34%  "
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "62.270.70"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "2.426.44"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").Value = "'11.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("D21").Value = "'325.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("D23").Value = "'2.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.23%  "
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").Value = "'65.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'620.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "'8.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("D28").Value = "0.0₃0961"
$ws.Range("E28").Value = "  -7.17%  "
$ws.Range("D29").Value = "2.548.77"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  -3.90%  "
$ws.Range("D32").Value = "'8.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.49%  "
$ws.Range("E33").Value = "  -1.89%  "
$ws.Range("E34").Value = "  -7.91%  "
$ws.Range("D35").Value = "'5.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("E37").Value = "  -4.73%  "
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").Value = "'18.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("D40").Value = "'146.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("D41").Value = "'5.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.25%  "
$ws.Range("E42").Value = "  -5.19%  "
$ws.Range("D43").Value = "'42.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "'2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.65%  "
$ws.Range("D46").Value = "'144.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("E48").Value = "  -3.98%  "
$ws.Range("D49").Value = "'0.593"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").Value = "'19.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.68%  "
$ws.Range("E51").Value = "  -3.42%  "
